$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 19; $i++) {
    if ($i -lt 10) {
        $suffix = "0$i"
    } else {
        $suffix = "$i"
    }
    $val = "06010101$suffix"
    $ws.Cells.Item($i + 1, 8).Value = $val
}

$ws.Range("H2:H20").Select()
